$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-28 and extend with new rows 29-35 per target dataset
# Columns: A=date(serial) B=first_median C=second_median D=open E=close F=high G=low

$ws.Cells.Item(2, 1).Value = 45831
$ws.Cells.Item(2, 2).Value = 761
$ws.Cells.Item(2, 3).Value = 1288
$ws.Cells.Item(2, 4).Value = 761
$ws.Cells.Item(2, 5).Value = 1288
$ws.Cells.Item(2, 6).Value = 1499
$ws.Cells.Item(2, 7).Value = 704

$ws.Cells.Item(3, 1).Value = 45832
$ws.Cells.Item(3, 2).Value = 769
$ws.Cells.Item(3, 3).Value = 1361
$ws.Cells.Item(3, 4).Value = 769
$ws.Cells.Item(3, 5).Value = 1361
$ws.Cells.Item(3, 6).Value = 1499
$ws.Cells.Item(3, 7).Value = 704

$ws.Cells.Item(4, 1).Value = 45833
$ws.Cells.Item(4, 2).Value = 737
$ws.Cells.Item(4, 3).Value = 1417
$ws.Cells.Item(4, 4).Value = 737
$ws.Cells.Item(4, 5).Value = 1417
$ws.Cells.Item(4, 6).Value = 1499
$ws.Cells.Item(4, 7).Value = 677

$ws.Cells.Item(5, 1).Value = 45834
$ws.Cells.Item(5, 2).Value = 733
$ws.Cells.Item(5, 3).Value = 1288
$ws.Cells.Item(5, 4).Value = 733
$ws.Cells.Item(5, 5).Value = 1288
$ws.Cells.Item(5, 6).Value = 1499
$ws.Cells.Item(5, 7).Value = 699

$ws.Cells.Item(6, 1).Value = 45835
$ws.Cells.Item(6, 2).Value = 733
$ws.Cells.Item(6, 3).Value = 1417
$ws.Cells.Item(6, 4).Value = 733
$ws.Cells.Item(6, 5).Value = 1417
$ws.Cells.Item(6, 6).Value = 1499
$ws.Cells.Item(6, 7).Value = 689

$ws.Cells.Item(7, 1).Value = 45836
$ws.Cells.Item(7, 2).Value = 733
$ws.Cells.Item(7, 3).Value = 1417
$ws.Cells.Item(7, 4).Value = 733
$ws.Cells.Item(7, 5).Value = 1417
$ws.Cells.Item(7, 6).Value = 1499
$ws.Cells.Item(7, 7).Value = 691

$ws.Cells.Item(8, 1).Value = 45838
$ws.Cells.Item(8, 2).Value = 752
$ws.Cells.Item(8, 3).Value = 1417
$ws.Cells.Item(8, 4).Value = 752
$ws.Cells.Item(8, 5).Value = 1417
$ws.Cells.Item(8, 6).Value = 1499
$ws.Cells.Item(8, 7).Value = 699

$ws.Cells.Item(9, 1).Value = 45839
$ws.Cells.Item(9, 2).Value = 759
$ws.Cells.Item(9, 3).Value = 1417
$ws.Cells.Item(9, 4).Value = 759
$ws.Cells.Item(9, 5).Value = 1417
$ws.Cells.Item(9, 6).Value = 1499
$ws.Cells.Item(9, 7).Value = 689

$ws.Cells.Item(10, 1).Value = 45841
$ws.Cells.Item(10, 2).Value = 747
$ws.Cells.Item(10, 3).Value = 1118.01
$ws.Cells.Item(10, 4).Value = 747
$ws.Cells.Item(10, 5).Value = 1118.01
$ws.Cells.Item(10, 6).Value = 1349
$ws.Cells.Item(10, 7).Value = 683

$ws.Cells.Item(11, 1).Value = 45842
$ws.Cells.Item(11, 2).Value = 975.44
$ws.Cells.Item(11, 3).Value = 1199
$ws.Cells.Item(11, 4).Value = 975.44
$ws.Cells.Item(11, 5).Value = 1199
$ws.Cells.Item(11, 6).Value = 1499
$ws.Cells.Item(11, 7).Value = 700

$ws.Cells.Item(12, 1).Value = 45845
$ws.Cells.Item(12, 2).Value = 760
$ws.Cells.Item(12, 3).Value = 1417
$ws.Cells.Item(12, 4).Value = 760
$ws.Cells.Item(12, 5).Value = 1417
$ws.Cells.Item(12, 6).Value = 1499
$ws.Cells.Item(12, 7).Value = 580

$ws.Cells.Item(13, 1).Value = 45846
$ws.Cells.Item(13, 2).Value = 760
$ws.Cells.Item(13, 3).Value = 1417
$ws.Cells.Item(13, 4).Value = 760
$ws.Cells.Item(13, 5).Value = 1417
$ws.Cells.Item(13, 6).Value = 1499
$ws.Cells.Item(13, 7).Value = 716

$ws.Cells.Item(14, 1).Value = 45847
$ws.Cells.Item(14, 2).Value = 760
$ws.Cells.Item(14, 3).Value = 850
$ws.Cells.Item(14, 4).Value = 760
$ws.Cells.Item(14, 5).Value = 850
$ws.Cells.Item(14, 6).Value = 1499
$ws.Cells.Item(14, 7).Value = 705

$ws.Cells.Item(15, 1).Value = 45848
$ws.Cells.Item(15, 2).Value = 741
$ws.Cells.Item(15, 3).Value = 852
$ws.Cells.Item(15, 4).Value = 741
$ws.Cells.Item(15, 5).Value = 852
$ws.Cells.Item(15, 6).Value = 1499
$ws.Cells.Item(15, 7).Value = 705

$ws.Cells.Item(16, 1).Value = 45849
$ws.Cells.Item(16, 2).Value = 738
$ws.Cells.Item(16, 3).Value = 1417
$ws.Cells.Item(16, 4).Value = 738
$ws.Cells.Item(16, 5).Value = 1417
$ws.Cells.Item(16, 6).Value = 1499
$ws.Cells.Item(16, 7).Value = 699

$ws.Cells.Item(17, 1).Value = 45850
$ws.Cells.Item(17, 2).Value = 733
$ws.Cells.Item(17, 3).Value = 1417
$ws.Cells.Item(17, 4).Value = 733
$ws.Cells.Item(17, 5).Value = 1417
$ws.Cells.Item(17, 6).Value = 1499
$ws.Cells.Item(17, 7).Value = 553

$ws.Cells.Item(18, 1).Value = 45851
$ws.Cells.Item(18, 2).Value = 733
$ws.Cells.Item(18, 3).Value = 1417
$ws.Cells.Item(18, 4).Value = 733
$ws.Cells.Item(18, 5).Value = 1417
$ws.Cells.Item(18, 6).Value = 1499
$ws.Cells.Item(18, 7).Value = 549

$ws.Cells.Item(19, 1).Value = 45852
$ws.Cells.Item(19, 2).Value = 737
$ws.Cells.Item(19, 3).Value = 1417
$ws.Cells.Item(19, 4).Value = 737
$ws.Cells.Item(19, 5).Value = 1417
$ws.Cells.Item(19, 6).Value = 1499
$ws.Cells.Item(19, 7).Value = 700

$ws.Cells.Item(20, 1).Value = 45853
$ws.Cells.Item(20, 2).Value = 737
$ws.Cells.Item(20, 3).Value = 1276
$ws.Cells.Item(20, 4).Value = 737
$ws.Cells.Item(20, 5).Value = 1276
$ws.Cells.Item(20, 6).Value = 1499
$ws.Cells.Item(20, 7).Value = 701

$ws.Cells.Item(21, 1).Value = 45854
$ws.Cells.Item(21, 2).Value = 744
$ws.Cells.Item(21, 3).Value = 1417
$ws.Cells.Item(21, 4).Value = 744
$ws.Cells.Item(21, 5).Value = 1417
$ws.Cells.Item(21, 6).Value = 1499
$ws.Cells.Item(21, 7).Value = 674

$ws.Cells.Item(22, 1).Value = 45855
$ws.Cells.Item(22, 2).Value = 730
$ws.Cells.Item(22, 3).Value = 1265
$ws.Cells.Item(22, 4).Value = 730
$ws.Cells.Item(22, 5).Value = 1265
$ws.Cells.Item(22, 6).Value = 1499
$ws.Cells.Item(22, 7).Value = 651

$ws.Cells.Item(23, 1).Value = 45856
$ws.Cells.Item(23, 2).Value = 721
$ws.Cells.Item(23, 3).Value = 1094
$ws.Cells.Item(23, 4).Value = 721
$ws.Cells.Item(23, 5).Value = 1094
$ws.Cells.Item(23, 6).Value = 1499
$ws.Cells.Item(23, 7).Value = 560

$ws.Cells.Item(24, 1).Value = 45857
$ws.Cells.Item(24, 2).Value = 752
$ws.Cells.Item(24, 3).Value = 1417
$ws.Cells.Item(24, 4).Value = 752
$ws.Cells.Item(24, 5).Value = 1417
$ws.Cells.Item(24, 6).Value = 1499
$ws.Cells.Item(24, 7).Value = 550

$ws.Cells.Item(25, 1).Value = 45859
$ws.Cells.Item(25, 2).Value = 719
$ws.Cells.Item(25, 3).Value = 1417
$ws.Cells.Item(25, 4).Value = 719
$ws.Cells.Item(25, 5).Value = 1417
$ws.Cells.Item(25, 6).Value = 1499
$ws.Cells.Item(25, 7).Value = 668

$ws.Cells.Item(26, 1).Value = 45860
$ws.Cells.Item(26, 2).Value = 718
$ws.Cells.Item(26, 3).Value = 1417
$ws.Cells.Item(26, 4).Value = 718
$ws.Cells.Item(26, 5).Value = 1417
$ws.Cells.Item(26, 6).Value = 1499
$ws.Cells.Item(26, 7).Value = 669

$ws.Cells.Item(27, 1).Value = 45861
$ws.Cells.Item(27, 2).Value = 739
$ws.Cells.Item(27, 3).Value = 1417
$ws.Cells.Item(27, 4).Value = 739
$ws.Cells.Item(27, 5).Value = 1417
$ws.Cells.Item(27, 6).Value = 1499
$ws.Cells.Item(27, 7).Value = 657

$ws.Cells.Item(28, 1).Value = 45862
$ws.Cells.Item(28, 2).Value = 726
$ws.Cells.Item(28, 3).Value = 1198.8
$ws.Cells.Item(28, 4).Value = 726
$ws.Cells.Item(28, 5).Value = 1198.8
$ws.Cells.Item(28, 6).Value = 1499
$ws.Cells.Item(28, 7).Value = 669

$ws.Cells.Item(29, 1).Value = 45863
$ws.Cells.Item(29, 2).Value = 728
$ws.Cells.Item(29, 3).Value = 1191
$ws.Cells.Item(29, 4).Value = 728
$ws.Cells.Item(29, 5).Value = 1191
$ws.Cells.Item(29, 6).Value = 1499
$ws.Cells.Item(29, 7).Value = 669

$ws.Cells.Item(30, 1).Value = 45864
$ws.Cells.Item(30, 2).Value = 929
$ws.Cells.Item(30, 3).Value = 1187.56
$ws.Cells.Item(30, 4).Value = 929
$ws.Cells.Item(30, 5).Value = 1187.56
$ws.Cells.Item(30, 6).Value = 1499
$ws.Cells.Item(30, 7).Value = 674

$ws.Cells.Item(31, 1).Value = 45865
$ws.Cells.Item(31, 2).Value = 725
$ws.Cells.Item(31, 3).Value = 1187.56
$ws.Cells.Item(31, 4).Value = 725
$ws.Cells.Item(31, 5).Value = 1187.56
$ws.Cells.Item(31, 6).Value = 1401
$ws.Cells.Item(31, 7).Value = 674

$ws.Cells.Item(32, 1).Value = 45866
$ws.Cells.Item(32, 2).Value = 1000
$ws.Cells.Item(32, 3).Value = 1100
$ws.Cells.Item(32, 4).Value = 1000
$ws.Cells.Item(32, 5).Value = 1100
$ws.Cells.Item(32, 6).Value = 1499
$ws.Cells.Item(32, 7).Value = 674

$ws.Cells.Item(33, 1).Value = 45867
$ws.Cells.Item(33, 2).Value = 1000
$ws.Cells.Item(33, 3).Value = 1187.56
$ws.Cells.Item(33, 4).Value = 1000
$ws.Cells.Item(33, 5).Value = 1187.56
$ws.Cells.Item(33, 6).Value = 1499
$ws.Cells.Item(33, 7).Value = 673

$ws.Cells.Item(34, 1).Value = 45868
$ws.Cells.Item(34, 2).Value = 969
$ws.Cells.Item(34, 3).Value = 949
$ws.Cells.Item(34, 4).Value = 969
$ws.Cells.Item(34, 5).Value = 949
$ws.Cells.Item(34, 6).Value = 10000
$ws.Cells.Item(34, 7).Value = 673

$ws.Cells.Item(35, 1).Value = 45869
$ws.Cells.Item(35, 2).Value = 729
$ws.Cells.Item(35, 3).Value = 1187.56
$ws.Cells.Item(35, 4).Value = 729
$ws.Cells.Item(35, 5).Value = 1187.56
$ws.Cells.Item(35, 6).Value = 1499
$ws.Cells.Item(35, 7).Value = 663

# Ensure newly added date cells (rows 29-35) use the same date number format as existing date column
for ($r = 29; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
